$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 05:59:00"
$ws.Cells.Item(3, 1).Value = "Total filas: 42"
$ws.Cells.Item(19, 1).Value = "05:59:00"
$ws.Cells.Item(19, 2).Value = "06:00"
$ws.Cells.Item(19, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(20, 2).Value = "06:04"
$ws.Cells.Item(20, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(20, 4).Value = 59
$ws.Cells.Item(21, 2).Value = "06:11"
$ws.Cells.Item(21, 3).Value = "215A_EL PATO"
$ws.Cells.Item(21, 4).Value = 66
$ws.Cells.Item(22, 1).Value = "05:59:00"
$ws.Cells.Item(22, 2).Value = "06:12"
$ws.Cells.Item(22, 3).Value = "215A_EL PATO"
$ws.Cells.Item(22, 4).Value = 13
$ws.Cells.Item(23, 1).Value = "05:59:00"
$ws.Cells.Item(23, 2).Value = "06:14"
$ws.Cells.Item(23, 3).Value = "225_HARAS DEL SUR"
$ws.Cells.Item(23, 4).Value = 15
$ws.Cells.Item(24, 1).Value = "05:59:00"
$ws.Cells.Item(24, 2).Value = "06:21"
$ws.Cells.Item(24, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(24, 4).Value = 22
$ws.Cells.Item(25, 1).Value = "05:59:00"
$ws.Cells.Item(25, 2).Value = "06:27"
$ws.Cells.Item(25, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(25, 4).Value = 28
$ws.Cells.Item(26, 1).Value = "04:38:41"
$ws.Cells.Item(26, 2).Value = "06:29"
$ws.Cells.Item(26, 3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(26, 4).Value = 111
$ws.Cells.Item(27, 1).Value = "05:59:00"
$ws.Cells.Item(27, 2).Value = "06:30"
$ws.Cells.Item(27, 3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(27, 4).Value = 31
$ws.Cells.Item(28, 2).Value = "06:31"
$ws.Cells.Item(28, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(28, 4).Value = 86
$ws.Cells.Item(29, 1).Value = "05:59:00"
$ws.Cells.Item(29, 2).Value = "06:32"
$ws.Cells.Item(29, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(29, 4).Value = 33
$ws.Cells.Item(30, 1).Value = "05:59:00"
$ws.Cells.Item(30, 2).Value = "06:44"
$ws.Cells.Item(30, 3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(30, 4).Value = 45
$ws.Cells.Item(30, 5).Value = "LP1912"
$ws.Cells.Item(31, 1).Value = "05:05:17"
$ws.Cells.Item(31, 2).Value = "06:46"
$ws.Cells.Item(31, 3).Value = "215C_EL PATO"
$ws.Cells.Item(31, 4).Value = 101
$ws.Cells.Item(31, 5).Value = "LP1912"
$ws.Cells.Item(32, 1).Value = "05:59:00"
$ws.Cells.Item(32, 2).Value = "06:47"
$ws.Cells.Item(32, 3).Value = "215C_EL PATO"
$ws.Cells.Item(32, 4).Value = 48
$ws.Cells.Item(32, 5).Value = "LP1912"
$ws.Cells.Item(33, 1).Value = "05:59:00"
$ws.Cells.Item(33, 2).Value = "07:00"
$ws.Cells.Item(33, 3).Value = "14_ABASTO"
$ws.Cells.Item(33, 4).Value = 61
$ws.Cells.Item(33, 5).Value = "LP1912"
$ws.Cells.Item(34, 1).Value = "05:59:00"
$ws.Cells.Item(34, 2).Value = "07:05"
$ws.Cells.Item(34, 3).Value = "15_ABASTO"
$ws.Cells.Item(34, 4).Value = 66
$ws.Cells.Item(34, 5).Value = "LP1912"
$ws.Cells.Item(35, 1).Value = "05:59:00"
$ws.Cells.Item(35, 2).Value = "07:05"
$ws.Cells.Item(35, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(35, 4).Value = 66
$ws.Cells.Item(35, 5).Value = "LP1912"
$ws.Cells.Item(36, 1).Value = "05:59:00"
$ws.Cells.Item(36, 2).Value = "07:07"
$ws.Cells.Item(36, 3).Value = "225_GOMEZ"
$ws.Cells.Item(36, 4).Value = 68
$ws.Cells.Item(36, 5).Value = "LP1912"
$ws.Cells.Item(37, 1).Value = "05:59:00"
$ws.Cells.Item(37, 2).Value = "07:12"
$ws.Cells.Item(37, 3).Value = "215A_EL PATO"
$ws.Cells.Item(37, 4).Value = 73
$ws.Cells.Item(37, 5).Value = "LP1912"
$ws.Cells.Item(38, 1).Value = "05:59:00"
$ws.Cells.Item(38, 2).Value = "07:16"
$ws.Cells.Item(38, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(38, 4).Value = 77
$ws.Cells.Item(38, 5).Value = "LP1912"
$ws.Cells.Item(39, 1).Value = "05:59:00"
$ws.Cells.Item(39, 2).Value = "07:21"
$ws.Cells.Item(39, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(39, 4).Value = 82
$ws.Cells.Item(39, 5).Value = "LP1912"
$ws.Cells.Item(40, 1).Value = "05:59:00"
$ws.Cells.Item(40, 2).Value = "07:23"
$ws.Cells.Item(40, 3).Value = "10_OLMOS"
$ws.Cells.Item(40, 4).Value = 84
$ws.Cells.Item(40, 5).Value = "LP1912"
$ws.Cells.Item(41, 1).Value = "05:59:00"
$ws.Cells.Item(41, 2).Value = "07:32"
$ws.Cells.Item(41, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(41, 4).Value = 93
$ws.Cells.Item(41, 5).Value = "LP1912"
$ws.Cells.Item(42, 1).Value = "05:59:00"
$ws.Cells.Item(42, 2).Value = "07:32"
$ws.Cells.Item(42, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(42, 4).Value = 93
$ws.Cells.Item(42, 5).Value = "LP1912"
$ws.Cells.Item(43, 1).Value = "05:59:00"
$ws.Cells.Item(43, 2).Value = "07:32"
$ws.Cells.Item(43, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(43, 4).Value = 93
$ws.Cells.Item(43, 5).Value = "LP1912"
$ws.Cells.Item(44, 1).Value = "05:59:00"
$ws.Cells.Item(44, 2).Value = "07:37"
$ws.Cells.Item(44, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(44, 4).Value = 98
$ws.Cells.Item(44, 5).Value = "LP1912"
$ws.Cells.Item(45, 1).Value = "05:59:00"
$ws.Cells.Item(45, 2).Value = "07:39"
$ws.Cells.Item(45, 3).Value = "10_OLMOS"
$ws.Cells.Item(45, 4).Value = 100
$ws.Cells.Item(45, 5).Value = "LP1912"
$ws.Cells.Item(46, 1).Value = "05:59:00"
$ws.Cells.Item(46, 2).Value = "07:48"
$ws.Cells.Item(46, 3).Value = "14_ABASTO"
$ws.Cells.Item(46, 4).Value = 109
$ws.Cells.Item(46, 5).Value = "LP1912"
$ws.Cells.Item(47, 1).Value = "05:59:00"
$ws.Cells.Item(47, 2).Value = "07:52"
$ws.Cells.Item(47, 3).Value = "215D_EL PATO"
$ws.Cells.Item(47, 4).Value = 113
$ws.Cells.Item(47, 5).Value = "LP1912"

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 05:59:00"
$ws.Cells.Item(3, 1).Value = "Total filas: 10"
$ws.Cells.Item(11, 1).Value = "05:59:00"
$ws.Cells.Item(11, 2).Value = "06:12"
$ws.Cells.Item(11, 3).Value = "215A_EL PATO"
$ws.Cells.Item(11, 4).Value = 13
$ws.Cells.Item(12, 1).Value = "05:05:17"
$ws.Cells.Item(12, 2).Value = "06:46"
$ws.Cells.Item(12, 3).Value = "215C_EL PATO"
$ws.Cells.Item(12, 4).Value = 101
$ws.Cells.Item(12, 5).Value = "LP1912"
$ws.Cells.Item(13, 1).Value = "05:59:00"
$ws.Cells.Item(13, 2).Value = "06:47"
$ws.Cells.Item(13, 3).Value = "215C_EL PATO"
$ws.Cells.Item(13, 4).Value = 48
$ws.Cells.Item(13, 5).Value = "LP1912"
$ws.Cells.Item(14, 1).Value = "05:59:00"
$ws.Cells.Item(14, 2).Value = "07:12"
$ws.Cells.Item(14, 3).Value = "215A_EL PATO"
$ws.Cells.Item(14, 4).Value = 73
$ws.Cells.Item(14, 5).Value = "LP1912"
$ws.Cells.Item(15, 1).Value = "05:59:00"
$ws.Cells.Item(15, 2).Value = "07:52"
$ws.Cells.Item(15, 3).Value = "215D_EL PATO"
$ws.Cells.Item(15, 4).Value = 113
$ws.Cells.Item(15, 5).Value = "LP1912"

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 05:59:00"
$ws.Cells.Item(3, 1).Value = "Total filas: 8"
$ws.Cells.Item(9, 1).Value = "05:59:00"
$ws.Cells.Item(9, 4).Value = 10
$ws.Cells.Item(11, 1).Value = "05:59:00"
$ws.Cells.Item(11, 4).Value = 34
$ws.Cells.Item(12, 1).Value = "05:59:00"
$ws.Cells.Item(12, 4).Value = 61
$ws.Cells.Item(13, 1).Value = "05:59:00"
$ws.Cells.Item(13, 2).Value = "07:35"
$ws.Cells.Item(13, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(13, 4).Value = 96
$ws.Cells.Item(13, 5).Value = "L6173"
